$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "64.118.25"
Set-TextValue "E2" "  -1.33%  "
Set-TextValue "D3" "3.521.32"
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "586.19"
Set-TextValue "E5" "  +0.15%  "
Set-TextValue "D6" "134.37"
Set-TextValue "E6" "  +0.41%  "
Set-TextValue "D7" "3.522.08"
Set-TextValue "E7" "  +0.21%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "E9" "  +0.20%  "
Set-TextValue "E10" "  -0.29%  "
Set-TextValue "D11" "7.12"
Set-TextValue "E11" "  -0.78%  "
Set-TextValue "D12" "0.377"
Set-TextValue "E12" "  -1.96%  "
Set-TextValue "D13" "4.123.95"
Set-TextValue "E13" "  +0.11%  "
Set-TextValue "D14" "27.48"
Set-TextValue "E14" "  -0.30%  "
Set-TextValue "E15" "  +1.43%  "
Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.530.52"
Set-TextValue "E16" "  +0.31%  "
Set-TextValue "B17" "ShibaInu"
Set-TextValue "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000178"
Set-TextValue "E17" "  -1.56%  "
Set-TextValue "D18" "64.183.13"
Set-TextValue "E18" "  -1.24%  "
Set-TextValue "D19" "9.79"
Set-TextValue "E19" "  -2.63%  "
Set-TextValue "D20" "13.87"
Set-TextValue "E20" "  -2.96%  "
Set-TextValue "D21" "5.61"
Set-TextValue "E21" "  -0.76%  "
Set-TextValue "D22" "383.12"
Set-TextValue "E22" "  -1.94%  "
Set-TextValue "B23" "WrappedeETH"
Set-TextValue "C23" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D23" "3.666.49"
Set-TextValue "E23" "  +0.19%  "
Set-TextValue "B24" "Polygon"
Set-TextValue "C24" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D24" "0.569"
Set-TextValue "E24" "  -1.08%  "
Set-TextValue "D25" "74.04"
Set-TextValue "E25" "  -0.96%  "
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "E27" "  -1.38%  "
Set-TextValue "E28" "  +3.47%  "
Set-TextValue "E29" "  -1.70%  "
Set-TextValue "E30" "  -1.50%  "
Set-TextValue "E31" "  -0.01%  "
Set-TextValue "D32" "8.45"
Set-TextValue "E32" "  +1.92%  "
Set-TextValue "E33" "  -1.09%  "
Set-TextValue "D34" "3.538.08"
Set-TextValue "E34" "  +0.40%  "
Set-TextValue "E35" "  -0.01%  "
Set-TextValue "D36" "23.58"
Set-TextValue "E36" "  -1.98%  "
Set-TextValue "D37" "0.145"
Set-TextValue "E37" "  -0.50%  "
Set-TextValue "D38" "5.39"
Set-TextValue "E38" "  +3.66%  "
Set-TextValue "D39" "6.94"
Set-TextValue "E39" "  +0.18%  "
Set-TextValue "E40" "  -0.80%  "
Set-TextValue "D41" "158.62"
Set-TextValue "E41" "  -6.36%  "
Set-TextValue "D42" "0.0787"
Set-TextValue "E42" "  -2.24%  "
Set-TextValue "D43" "26.55"
Set-TextValue "E43" "  +1.45%  "
Set-TextValue "E44" "  -0.82%  "
Set-TextValue "E45" "  +0.06%  "
Set-TextValue "E46" "  -2.67%  "
Set-TextValue "D47" "41.62"
Set-TextValue "E47" "  -3.11%  "
Set-TextValue "D48" "4.39"
Set-TextValue "E48" "  -0.56%  "
Set-TextValue "E49" "  -2.35%  "
Set-TextValue "D50" "2.481.52"
Set-TextValue "E50" "  +0.35%  "
Set-TextValue "D51" "6.80"
Set-TextValue "E51" "  -0.86%  "
